# Updated TestData for Portugal Market
#
# - Adds a new "Portugal" worksheet (cloned from "Swiss") right after "Swiss"
# - Populates it with the Portugal market name / Jira reference
# - Makes "Portugal" the active sheet/tab with B4 selected
# - Leaves "Swiss" selected over its full used range (A1:D11), no longer the active tab

$wb = $excel.ActiveWorkbook

# Change the selection on the Swiss sheet to the whole used range before duplicating it,
# so the clone inherits that sheetView and we don't leave the old B2:B4 selection behind.
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Activate()
$swiss.Range("A1:D11").Select()

# Duplicate "Swiss" to create the new sheet, inserting the copy immediately after it.
$swiss.Copy([System.Reflection.Missing]::Value, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Fill in the Portugal-specific market name and Jira reference.
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2404"

# Rows 3 & 4 grow to a taller (two-line) height on the new sheet.
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8

# Make the new sheet the active tab, with B4 selected.
$portugal.Activate()
$portugal.Range("B4").Select()
